# Update "想去人数" (F column) counts for "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(5, 6).Value = 300
$ws1.Cells.Item(6, 6).Value = 1122
$ws1.Cells.Item(7, 6).Value = 1465
$ws1.Cells.Item(8, 6).Value = 588
$ws1.Cells.Item(10, 6).Value = 766
$ws1.Cells.Item(11, 6).Value = 72
$ws1.Cells.Item(14, 6).Value = 465
$ws1.Cells.Item(15, 6).Value = 1411
$ws1.Cells.Item(16, 6).Value = 133
$ws1.Cells.Item(20, 6).Value = 89
$ws1.Cells.Item(21, 6).Value = 667
$ws1.Cells.Item(22, 6).Value = 1022
$ws1.Cells.Item(24, 6).Value = 258
$ws1.Cells.Item(25, 6).Value = 32
$ws1.Cells.Item(26, 6).Value = 6076
$ws1.Cells.Item(31, 6).Value = 14799
$ws1.Cells.Item(34, 6).Value = 112
$ws1.Cells.Item(35, 6).Value = 91
$ws1.Cells.Item(36, 6).Value = 10825
$ws1.Cells.Item(37, 6).Value = 671
$ws1.Cells.Item(38, 6).Value = 4237
$ws1.Cells.Item(39, 6).Value = 175

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 300
$ws4.Cells.Item(6, 6).Value = 1122
$ws4.Cells.Item(7, 6).Value = 1465
$ws4.Cells.Item(8, 6).Value = 588
$ws4.Cells.Item(10, 6).Value = 766
$ws4.Cells.Item(11, 6).Value = 72
$ws4.Cells.Item(14, 6).Value = 465
$ws4.Cells.Item(15, 6).Value = 1411
$ws4.Cells.Item(16, 6).Value = 133
$ws4.Cells.Item(21, 6).Value = 89
$ws4.Cells.Item(22, 6).Value = 667
$ws4.Cells.Item(24, 6).Value = 1022
$ws4.Cells.Item(26, 6).Value = 258
$ws4.Cells.Item(27, 6).Value = 32
$ws4.Cells.Item(29, 6).Value = 6076
$ws4.Cells.Item(34, 6).Value = 14799
$ws4.Cells.Item(37, 6).Value = 112
$ws4.Cells.Item(38, 6).Value = 91
$ws4.Cells.Item(39, 6).Value = 10825
$ws4.Cells.Item(40, 6).Value = 671
$ws4.Cells.Item(41, 6).Value = 4237
$ws4.Cells.Item(42, 6).Value = 175

